# Server: disconnect한 player delete 시 오류 수정
# Update the "Move" sheet's Value column for Ch_Walk (20002) and Ch_Run (20003)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Move")

# C3 = Ch_Walk Value: 0.5 -> 1
$ws.Range("C3").Value = 1

# C4 = Ch_Run Value: 0.8 -> 1.3
$ws.Range("C4").Value = 1.3
